## ---------------------------------------------------------------------
## parent-template.xlsx edit: add religion/gender/gsm/username columns,
## data validation, a NOTE legend sheet, and rename Sheet1 -> Sheet2.
## ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Rename the original (only) sheet to "Sheet2" ------------------
$data = $wb.Worksheets.Item(1)
$data.Name = "Sheet2"

# ---- 2. Header row: write in the same order the original author must
#        have used so shared-string indices line up with the target.
#        (A,B,C,F,G,H,I first, then D,E last)
$data.Range("A1").Value = "firstname"
$data.Range("B1").Value = "surname"
$data.Range("C1").Value = "othername"
$data.Range("F1").Value = "email"
$data.Range("G1").Value = "address"
$data.Range("H1").Value = "religion"
$data.Range("I1").Value = "gender"
$data.Range("D1").Value = "gsm"
$data.Range("E1").Value = "username"

# Red header text on firstname / surname / gsm / address
$data.Range("A1").Font.Color = 255
$data.Range("B1").Font.Color = 255
$data.Range("D1").Font.Color = 255
$data.Range("G1").Font.Color = 255

# ---- 3. Data validations ----------------------------------------------
$vI = $data.Range("I1:I1048576").Validation
$vI.Add(1, 1, 1, 1, 3)
$vI.ErrorMessage = "use 1 for Male, 2 for Female and 3 for other"

$vH = $data.Range("H1:H1048576").Validation
$vH.Add(1, 1, 1, 1, 3)
$vH.ErrorMessage = "use 1 for Muslim, 2 for Christain and 3 for Other"

$vD = $data.Range("D1:D1048576").Validation
$vD.Add(6, 1, 3, 10)
$vD.IgnoreBlank = $false
$vD.ErrorMessage = "10 digits, omit the first 0"
$vD.InputTitle = "phone number"

$vB = $data.Range("B1:B1048576").Validation
$vB.Add(6, 1, 6, 30)
$vB.IgnoreBlank = $false
$vB.ErrorMessage = "maximum of 30 character"

$vA = $data.Range("A1:A1048576").Validation
$vA.Add(6, 1, 6, 31)
$vA.IgnoreBlank = $true
$vA.ErrorMessage = "maximum of 30 character"

$vG = $data.Range("G1").Validation
$vG.Add(6, 1, 6, 256)
$vG.IgnoreBlank = $true
$vG.ErrorMessage = "address is compulsary and not more than 255 characters"

# ---- 4. Column widths (bestFit-style widths from the template) --------
$data.Columns.Item(1).ColumnWidth = 8.666666666666666
$data.Columns.Item(2).ColumnWidth = 7.833333333333333
$data.Columns.Item(3).ColumnWidth = 10
$data.Columns.Item(4).ColumnWidth = 10.166666666666666
$data.Columns.Item(5).ColumnWidth = 9

# ---- 5. Freeze header row, set selection -------------------------------
$data.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$data.Range("H8").Select()

# ---- 6. Add the NOTE sheet after Sheet2 --------------------------------
$note = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$note.Name = "NOTE"

# Row 1: GENDER / RELIGION banners
$note.Range("A1").Value = "GENDER "
$note.Range("D1").Value = "RELIGION"
$note.Range("A1:B1").Font.Size = 18
$note.Range("A1:B1").Font.Color = 255
$note.Range("A1:B1").HorizontalAlignment = -4108
$note.Range("A1:B1").WrapText = $true
$note.Range("A1:B1").Merge()
$note.Range("D1:E1").Font.Size = 16
$note.Range("D1:E1").Font.Color = 255
$note.Range("D1:E1").HorizontalAlignment = -4108
$note.Range("D1:E1").Merge()
$note.Rows.Item(1).RowHeight = 23.25

# Row 2: explanation text
$note.Range("A2").Value = "the number coresponding number represent each gender"
$note.Range("D2").Value = "the corresponding number represent each religion"
$note.Range("A2:B2").Font.Size = 18
$note.Range("A2:B2").HorizontalAlignment = -4108
$note.Range("A2:B2").WrapText = $true
$note.Range("A2:B2").Borders.Item(9).LineStyle = 1
$note.Range("A2:B2").Merge()
$note.Range("D2:E2").Font.Size = 16
$note.Range("D2:E2").HorizontalAlignment = -4108
$note.Range("D2:E2").WrapText = $true
$note.Range("D2:E2").Borders.Item(9).LineStyle = 1
$note.Range("D2:E2").Merge()
$note.Rows.Item(2).RowHeight = 23.25

# Rows 3-5: legend values
$note.Range("A3").Value = "MALE"
$note.Range("B3").Value = 1
$note.Range("D3").Value = "MUSLIM"
$note.Range("E3").Value = 1

$note.Range("A4").Value = "FEMALE"
$note.Range("B4").Value = 2
$note.Range("D4").Value = "CHRISTIAN"
$note.Range("E4").Value = 2

$note.Range("A5").Value = "OTHER"
$note.Range("B5").Value = 3
$note.Range("D5").Value = "OTHER"
$note.Range("E5").Value = 3

$note.Range("A3:A5").Font.Size = 48
$note.Range("A3:A5").Font.Color = 16777215
$note.Range("A3:A5").Interior.Color = 5287936
$note.Range("A3:A5").Borders.LineStyle = 1

$note.Range("B3:B5").Font.Size = 48
$note.Range("B3:B5").Font.Color = 16777215
$note.Range("B3:B5").Interior.Color = 5287936
$note.Range("B3:B5").Borders.LineStyle = 1
$note.Range("B3:B5").HorizontalAlignment = -4108

$note.Range("D3:D5").Font.Size = 48
$note.Range("D3:D5").Font.Color = 16777215
$note.Range("D3:D5").Interior.Color = 192
$note.Range("D3:D5").Borders.LineStyle = 1

$note.Range("E3:E5").Font.Size = 48
$note.Range("E3:E5").Font.Color = 16777215
$note.Range("E3:E5").Interior.Color = 192
$note.Range("E3:E5").Borders.LineStyle = 1
$note.Range("E3:E5").HorizontalAlignment = -4108
$note.Range("E3:E5").VerticalAlignment = -4108

$note.Rows.Item(3).RowHeight = 61.5
$note.Rows.Item(4).RowHeight = 61.5
$note.Rows.Item(5).RowHeight = 61.5

# Row 8: leftover formatted placeholder cells
$note.Range("D8:E8").Font.Size = 18
$note.Range("D8:E8").Font.Color = 255
$note.Range("D8:E8").Borders.Item(9).LineStyle = 1
$note.Range("D8:E8").Merge()
$note.Rows.Item(8).RowHeight = 23.25

# Column widths for NOTE sheet
$note.Columns.Item(1).ColumnWidth = 32.666666666666664
$note.Columns.Item(2).ColumnWidth = 17.28515625
$note.Columns.Item(4).ColumnWidth = 42.666666666666664
$note.Columns.Item(5).ColumnWidth = 10.333333333333334

$note.Range("B13").Select()

# ---- 7. Leave Sheet2 as the active/selected sheet ----------------------
$data.Activate()

Write-Output "done"
